# Generate Report for Handback
# Populates the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns on the zh-cn and de-de sheets now that handback has completed, updates the
# "Status" column text/column-width on all three sheets, and widens a couple of columns
# that now hold longer content.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: status text + column widths
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# Helper: stamp the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns (I, J, K) for a locale sheet, for both
# data rows (2 and 3), and widen the columns that now hold real content.
# ---------------------------------------------------------------------------
function Set-HandbackInfo($sheet, $mdUrl, $xlfFileName, $handbackDateTime) {

    $sheet.Columns.Item(3).ColumnWidth = 29.166666666666668
    $sheet.Columns.Item(9).ColumnWidth = 39.166666666666664
    $sheet.Columns.Item(10).ColumnWidth = 39.166666666666664

    foreach ($row in 2, 3) {
        $iCell = $sheet.Cells.Item($row, 9)
        $iCell.Value = "4623fd6e-19e1-49af-b8f1-6fc6bc73d7a3.md"
        $sheet.Hyperlinks.Add($iCell, $mdUrl, "", "", "4623fd6e-19e1-49af-b8f1-6fc6bc73d7a3.md")
        $iCell.Font.Underline = $true
        $iCell.Font.Color = 15570276

        $sheet.Cells.Item($row, 10).Value = $xlfFileName
        $sheet.Cells.Item($row, 11).Value = $handbackDateTime
    }
}

$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9239d90d3bd53e1f6ac10a6f73a9e3a8cb272d43/e2e/4623fd6e-19e1-49af-b8f1-6fc6bc73d7a3.md"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
Set-HandbackInfo $zhcn $mdUrl "4623fd6e-19e1-49af-b8f1-6fc6bc73d7a3.b4fb12fd6d03f41cd78bb575a45f29b024da6344.zh-cn.xlf" "2016-08-17 17:01:55"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
Set-HandbackInfo $dede $mdUrl "4623fd6e-19e1-49af-b8f1-6fc6bc73d7a3.b4fb12fd6d03f41cd78bb575a45f29b024da6344.de-de.xlf" "2016-08-17 17:02:09"
